$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 2.728779965065642
$ws.Cells.Item(2, 4).Value = 9.95505111298154
$ws.Cells.Item(2, 5).Value = 13.89494703328383
$ws.Cells.Item(2, 6).Value = 29.25823648110417
$ws.Cells.Item(2, 7).Value = 30.1911326577861
$ws.Cells.Item(2, 8).Value = 13.72557788070845
$ws.Cells.Item(2, 9).Value = 19.95349761985225
$ws.Cells.Item(2, 10).Value = 9.686377412883518
$ws.Cells.Item(2, 14).Value = 18.99769123883711
$ws.Cells.Item(2, 15).Value = 21.27759823122737
$ws.Cells.Item(3, 3).Value = 2.666799860447807
$ws.Cells.Item(3, 4).Value = 9.954949846502856
$ws.Cells.Item(3, 5).Value = 13.86112943064488
$ws.Cells.Item(3, 6).Value = 28.95767149847613
$ws.Cells.Item(3, 7).Value = 29.49286781126834
$ws.Cells.Item(3, 8).Value = 13.68904849605218
$ws.Cells.Item(3, 9).Value = 19.77431705042589
$ws.Cells.Item(3, 10).Value = 9.688278086992334
$ws.Cells.Item(3, 14).Value = 18.40031508502701
$ws.Cells.Item(3, 15).Value = 21.09485627290016
$ws.Cells.Item(4, 3).Value = 2.627431768047288
$ws.Cells.Item(4, 4).Value = 9.956579151903831
$ws.Cells.Item(4, 5).Value = 13.84336193971557
$ws.Cells.Item(4, 6).Value = 28.77966669272912
$ws.Cells.Item(4, 7).Value = 29.06527682824506
$ws.Cells.Item(4, 8).Value = 13.66951148769116
$ws.Cells.Item(4, 9).Value = 19.66886519832054
$ws.Cells.Item(4, 10).Value = 9.691131607468559
$ws.Cells.Item(4, 14).Value = 18.02485520896362
$ws.Cells.Item(4, 15).Value = 20.98749649895068
$ws.Cells.Item(5, 3).Value = 2.611069552182935
$ws.Cells.Item(5, 4).Value = 9.957669060315016
$ws.Cells.Item(5, 5).Value = 13.83687981959222
$ws.Cells.Item(5, 6).Value = 28.70885334558064
$ws.Cells.Item(5, 7).Value = 28.89160399499856
$ws.Cells.Item(5, 8).Value = 13.66228266299906
$ws.Cells.Item(5, 9).Value = 19.62708502643778
$ws.Cells.Item(5, 10).Value = 9.692718850239771
$ws.Cells.Item(5, 14).Value = 17.86990355188765
$ws.Cells.Item(5, 15).Value = 20.94500837609291
$ws.Cells.Item(6, 3).Value = 2.608333649158011
$ws.Cells.Item(6, 4).Value = 9.957875779935625
$ws.Cells.Item(6, 5).Value = 13.83584939071642
$ws.Cells.Item(6, 6).Value = 28.69720128737304
$ws.Cells.Item(6, 7).Value = 28.86280889360314
$ws.Cells.Item(6, 8).Value = 13.66112671683065
$ws.Cells.Item(6, 9).Value = 19.62022071728637
$ws.Cells.Item(6, 10).Value = 9.693008050396317
$ws.Cells.Item(6, 14).Value = 17.84406337566597
$ws.Cells.Item(6, 15).Value = 20.93803069423352
$ws.Cells.Item(7, 3).Value = 2.627212380297049
$ws.Cells.Item(7, 4).Value = 9.956592125435105
$ws.Cells.Item(7, 5).Value = 13.84327144362124
$ws.Cells.Item(7, 6).Value = 28.77870459204133
$ws.Cells.Item(7, 7).Value = 29.06293190951466
$ws.Cells.Item(7, 8).Value = 13.66941102414877
$ws.Cells.Item(7, 9).Value = 19.66829685213405
$ws.Cells.Item(7, 10).Value = 9.691151294850247
$ws.Cells.Item(7, 14).Value = 18.02277304767603
$ws.Cells.Item(7, 15).Value = 20.9869183243332
$ws.Cells.Item(8, 3).Value = 2.707687640181533
$ws.Cells.Item(8, 4).Value = 9.954665546593931
$ws.Cells.Item(8, 5).Value = 13.88266783122267
$ws.Cells.Item(8, 6).Value = 29.15329252787406
$ws.Cells.Item(8, 7).Value = 29.950302168238
$ws.Cells.Item(8, 8).Value = 13.71238573299217
$ws.Cells.Item(8, 9).Value = 19.89079580457227
$ws.Cells.Item(8, 10).Value = 9.686682972881471
$ws.Cells.Item(8, 14).Value = 18.79364780656867
$ws.Cells.Item(8, 15).Value = 21.21361027950495
$ws.Cells.Item(9, 3).Value = 2.854662956744316
$ws.Cells.Item(9, 4).Value = 9.964275670941229
$ws.Cells.Item(9, 5).Value = 13.98346104558476
$ws.Cells.Item(9, 6).Value = 29.9361084707968
$ws.Cells.Item(9, 7).Value = 31.68738879936571
$ws.Cells.Item(9, 8).Value = 13.8193383253815
$ws.Cells.Item(9, 9).Value = 20.36128806603421
$ws.Cells.Item(9, 10).Value = 9.69128245322317
$ws.Cells.Item(9, 14).Value = 20.2273683202997
$ws.Cells.Item(9, 15).Value = 21.6945572092241
$ws.Cells.Item(10, 3).Value = 2.955520451787607
$ws.Cells.Item(10, 4).Value = 9.979444865072566
$ws.Cells.Item(10, 5).Value = 14.07149579450354
$ws.Cells.Item(10, 6).Value = 30.53555718835969
$ws.Cells.Item(10, 7).Value = 32.94627397956579
$ws.Cells.Item(10, 8).Value = 13.91133742855725
$ws.Cells.Item(10, 9).Value = 20.72490527758725
$ws.Cells.Item(10, 10).Value = 9.702772822597744
$ws.Cells.Item(10, 14).Value = 21.22223697909767
$ws.Cells.Item(10, 15).Value = 22.06723566228217
$ws.Cells.Item(11, 3).Value = 2.999759994391921
$ws.Cells.Item(11, 4).Value = 9.988091451713284
$ws.Cells.Item(11, 5).Value = 14.11448837401471
$ws.Cells.Item(11, 6).Value = 30.81242098004964
$ws.Cells.Item(11, 7).Value = 33.51218979878003
$ws.Cells.Item(11, 8).Value = 13.95600390532045
$ws.Cells.Item(11, 9).Value = 20.89357776281312
$ws.Cells.Item(11, 10).Value = 9.709750526752845
$ws.Cells.Item(11, 14).Value = 21.66018057919901
$ws.Cells.Item(11, 15).Value = 22.24033131613683
$ws.Cells.Item(12, 3).Value = 3.016269201105819
$ws.Cells.Item(12, 4).Value = 9.99161515712877
$ws.Cells.Item(12, 5).Value = 14.13118312209874
$ws.Cells.Item(12, 6).Value = 30.91776706014782
$ws.Cells.Item(12, 7).Value = 33.7252850667158
$ws.Cells.Item(12, 8).Value = 13.97331327431802
$ws.Cells.Item(12, 9).Value = 20.95786339543769
$ws.Cells.Item(12, 10).Value = 9.712643302057842
$ws.Cells.Item(12, 14).Value = 21.82377585682186
$ws.Cells.Item(12, 15).Value = 22.30633505977262
$ws.Cells.Item(13, 3).Value = 3.012724581529373
$ws.Cells.Item(13, 4).Value = 9.990845198871178
$ws.Cells.Item(13, 5).Value = 14.12756933111594
$ws.Cells.Item(13, 6).Value = 30.89505801903289
$ws.Cells.Item(13, 7).Value = 33.67944830912114
$ws.Cells.Item(13, 8).Value = 13.96956797943543
$ws.Cells.Item(13, 9).Value = 20.94400085894364
$ws.Cells.Item(13, 10).Value = 9.712009173622972
$ws.Cells.Item(13, 14).Value = 21.78864458690801
$ws.Cells.Item(13, 15).Value = 22.29210058217868
$ws.Cells.Item(14, 3).Value = 3.001123133745236
$ws.Cells.Item(14, 4).Value = 9.988376360225654
$ws.Cells.Item(14, 5).Value = 14.11585361043715
$ws.Cells.Item(14, 6).Value = 30.82107831877407
$ws.Cells.Item(14, 7).Value = 33.52974669752561
$ws.Cells.Item(14, 8).Value = 13.95742009621922
$ws.Cells.Item(14, 9).Value = 20.8988586244736
$ws.Cells.Item(14, 10).Value = 9.709983505264406
$ws.Cells.Item(14, 14).Value = 21.67368539489659
$ws.Cells.Item(14, 15).Value = 22.24575266440345
$ws.Cells.Item(15, 3).Value = 2.99398501318458
$ws.Cells.Item(15, 4).Value = 9.98689655637731
$ws.Cells.Item(15, 5).Value = 14.10873107781528
$ws.Cells.Item(15, 6).Value = 30.77582637344956
$ws.Cells.Item(15, 7).Value = 33.43788660025945
$ws.Cells.Item(15, 8).Value = 13.9500303268308
$ws.Cells.Item(15, 9).Value = 20.8712598241902
$ws.Cells.Item(15, 10).Value = 9.708775302607048
$ws.Cells.Item(15, 14).Value = 21.60297336126124
$ws.Cells.Item(15, 15).Value = 22.21742096087996
$ws.Cells.Item(16, 3).Value = 2.952595627114936
$ws.Cells.Item(16, 4).Value = 9.978914788577089
$ws.Cells.Item(16, 5).Value = 14.0687446273984
$ws.Cells.Item(16, 6).Value = 30.51753921407916
$ws.Cells.Item(16, 7).Value = 32.90913378270817
$ws.Cells.Item(16, 8).Value = 13.90847420835358
$ws.Cells.Item(16, 9).Value = 20.71394298549198
$ws.Cells.Item(16, 10).Value = 9.702351936432436
$ws.Cells.Item(16, 14).Value = 21.19330956972086
$ws.Cells.Item(16, 15).Value = 22.05599035625912
$ws.Cells.Item(17, 3).Value = 2.926778821510825
$ws.Cells.Item(17, 4).Value = 9.974464323503957
$ws.Cells.Item(17, 5).Value = 14.04496227410564
$ws.Cells.Item(17, 6).Value = 30.36008846342321
$ws.Cells.Item(17, 7).Value = 32.582859400511
$ws.Cells.Item(17, 8).Value = 13.88369529259143
$ws.Cells.Item(17, 9).Value = 20.61822981597834
$ws.Cells.Item(17, 10).Value = 9.698858949691031
$ws.Cells.Item(17, 14).Value = 20.9381421901699
$ws.Cells.Item(17, 15).Value = 21.95783059923031
$ws.Cells.Item(18, 3).Value = 2.911775614143865
$ws.Cells.Item(18, 4).Value = 9.972068990275975
$ws.Cells.Item(18, 5).Value = 14.03156109064439
$ws.Cells.Item(18, 6).Value = 30.26992536866286
$ws.Cells.Item(18, 7).Value = 32.39457190214168
$ws.Cells.Item(18, 8).Value = 13.86970875867027
$ws.Cells.Item(18, 9).Value = 20.56348847637524
$ws.Cells.Item(18, 10).Value = 9.697014745871316
$ws.Cells.Item(18, 14).Value = 20.79000725568364
$ws.Cells.Item(18, 15).Value = 21.90171041309424
$ws.Cells.Item(19, 3).Value = 2.906669560475217
$ws.Cells.Item(19, 4).Value = 9.971286264311066
$ws.Cells.Item(19, 5).Value = 14.02707165926451
$ws.Cells.Item(19, 6).Value = 30.23946905090525
$ws.Cells.Item(19, 7).Value = 32.33072133737637
$ws.Cells.Item(19, 8).Value = 13.8650190696841
$ws.Cells.Item(19, 9).Value = 20.54500899855116
$ws.Cells.Item(19, 10).Value = 9.69641868259675
$ws.Cells.Item(19, 14).Value = 20.73962067985786
$ws.Cells.Item(19, 15).Value = 21.88276895557809
$ws.Cells.Item(20, 3).Value = 2.92954306496494
$ws.Cells.Item(20, 4).Value = 9.974921073813031
$ws.Cells.Item(20, 5).Value = 14.04746525924713
$ws.Cells.Item(20, 6).Value = 30.37680882686276
$ws.Cells.Item(20, 7).Value = 32.61765809388397
$ws.Cells.Item(20, 8).Value = 13.88630562597045
$ws.Cells.Item(20, 9).Value = 20.62838693071205
$ws.Cells.Item(20, 10).Value = 9.699213729356158
$ws.Cells.Item(20, 14).Value = 20.96544799484618
$ws.Cells.Item(20, 15).Value = 21.96824522223603
$ws.Cells.Item(21, 3).Value = 3.004537420118898
$ws.Cells.Item(21, 4).Value = 9.989094762572753
$ws.Cells.Item(21, 5).Value = 14.11928363092784
$ws.Cells.Item(21, 6).Value = 30.84279504255863
$ws.Cells.Item(21, 7).Value = 33.57375216391713
$ws.Cells.Item(21, 8).Value = 13.96097758055038
$ws.Cells.Item(21, 9).Value = 20.91210723109223
$ws.Cells.Item(21, 10).Value = 9.710571706527952
$ws.Cells.Item(21, 14).Value = 21.70751365554064
$ws.Cells.Item(21, 15).Value = 22.25935424026509
$ws.Cells.Item(22, 3).Value = 3.052129177597336
$ws.Cells.Item(22, 4).Value = 9.999811063477832
$ws.Cells.Item(22, 5).Value = 14.16863180903648
$ws.Cells.Item(22, 6).Value = 31.15023604654358
$ws.Cells.Item(22, 7).Value = 34.19149596608521
$ws.Cells.Item(22, 8).Value = 14.01207827847189
$ws.Cells.Item(22, 9).Value = 21.09991615645727
$ws.Cells.Item(22, 10).Value = 9.719454015011911
$ws.Cells.Item(22, 14).Value = 22.17935961385674
$ws.Cells.Item(22, 15).Value = 22.45224346823413
$ws.Cells.Item(23, 3).Value = 3.026860887910664
$ws.Cells.Item(23, 4).Value = 9.993959217700151
$ws.Cells.Item(23, 5).Value = 14.14207635031675
$ws.Cells.Item(23, 6).Value = 30.98591603218996
$ws.Cells.Item(23, 7).Value = 33.86251865599016
$ws.Cells.Item(23, 8).Value = 13.98459795704715
$ws.Cells.Item(23, 9).Value = 20.9994796571583
$ws.Cells.Item(23, 10).Value = 9.714580309249394
$ws.Cells.Item(23, 14).Value = 21.92877110912574
$ws.Cells.Item(23, 15).Value = 22.34907260818606
$ws.Cells.Item(24, 3).Value = 2.928293850742859
$ws.Cells.Item(24, 4).Value = 9.974714068218688
$ws.Cells.Item(24, 5).Value = 14.04633281295986
$ws.Cells.Item(24, 6).Value = 30.36924842973586
$ws.Cells.Item(24, 7).Value = 32.60192779189882
$ws.Cells.Item(24, 8).Value = 13.88512468607979
$ws.Cells.Item(24, 9).Value = 20.62379400570166
$ws.Cells.Item(24, 10).Value = 9.699052822652378
$ws.Cells.Item(24, 14).Value = 20.95310750188672
$ws.Cells.Item(24, 15).Value = 21.96353579133695
$ws.Cells.Item(25, 3).Value = 2.816118648598807
$ws.Cells.Item(25, 4).Value = 9.960247339920642
$ws.Cells.Item(25, 5).Value = 13.95370903001434
$ws.Cells.Item(25, 6).Value = 29.71968912420405
$ws.Cells.Item(25, 7).Value = 31.21943570676743
$ws.Cells.Item(25, 8).Value = 13.78801584276746
$ws.Cells.Item(25, 9).Value = 20.23064075146675
$ws.Cells.Item(25, 10).Value = 9.688610579383477
$ws.Cells.Item(25, 14).Value = 19.84905939529497
$ws.Cells.Item(25, 15).Value = 21.56084013990475
